# Adds a new weekly price record for "Alcachofa" (Vega Modelo de Temuco)
# at the top of the data block: a new row is inserted at row 123, pushing
# the existing rows 123-148 down to 124-149, and the new row 123 is filled
# with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 123 - shifts rows 123:148 down to 124:149
$ws.Rows("123:123").Insert()

# Populate the newly inserted row 123 with the new record
$ws.Cells.Item(123, 1).Value  = 10
$ws.Cells.Item(123, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(123, 3).Value  = 'La Araucanía'
$ws.Cells.Item(123, 4).Value  = 44508
$ws.Cells.Item(123, 5).Value  = 9
$ws.Cells.Item(123, 6).Value  = 100112013
$ws.Cells.Item(123, 7).Value  = 'Alcachofa'
$ws.Cells.Item(123, 8).Value  = 'Española'
$ws.Cells.Item(123, 9).Value  = 'Primera'
$ws.Cells.Item(123, 10).Value = 3000
$ws.Cells.Item(123, 11).Value = 400
$ws.Cells.Item(123, 12).Value = 400
$ws.Cells.Item(123, 13).Value = 400
$ws.Cells.Item(123, 14).Value = '$/unidad'
$ws.Cells.Item(123, 15).Value = 'Región del Maule'
$ws.Cells.Item(123, 16).Value = 400
$ws.Cells.Item(123, 17).Value = 1
$ws.Cells.Item(123, 18).Value = 'Hortaliza'
